$d = $word.ActiveDocument

# --- Helper: insert a brand-new empty paragraph right after $afterPara, return it ---
function New-ParaAfter($afterPara) {
    $r = $afterPara.Range
    $r.Collapse(0)
    $r.InsertParagraphAfter()
    return $afterPara.Next()
}

# Anchor paragraphs in the original document:
#   3: "Labb1.1" (bold heading)
#   4: "Jag hade inga speciella problem..."
#   5: empty paragraph
#   6: paragraph holding the _GoBack bookmark (currently empty)
$pEmpty5 = $d.Paragraphs(5)

# --- Insert the new "Labb1.2" block + 4 body paragraphs, all BEFORE the bookmark
#     paragraph. Create every paragraph first (plain, unformatted) and only fill
#     in text -- bold is applied at the very end so it doesn't leak into later
#     inserted paragraphs (new paragraph marks inherit formatting from the mark
#     they were split off from). ---
$pLabb12   = New-ParaAfter $pEmpty5
$pStarted  = New-ParaAfter $pLabb12
$pGenererade100x3 = New-ParaAfter $pStarted
$pGenererade100x10 = New-ParaAfter $pGenererade100x3
$pTestade  = New-ParaAfter $pGenererade100x10

$pLabb12.Range.Text = "Labb1.2"
$pStarted.Range.Text = "Började med att göra ett problem med 9 lådor och 9 skadade personer, det tog 0,37 sekunder att köra det problemet med IPP och tänkte då att jag måste göra problemen bra mycket större."
$pGenererade100x3.Range.Text = "Genererade då problem med 100 lådor samt 3 stycken skadade personer, även det gick ganska fort att köra med IPP, någon sekund bara."
$pGenererade100x10.Range.Text = "Genererade då ett problem med 100 lådor och 10 stycken skadade personer, det gick inte alls att köra med IPP och det tog även ett tag att köra med FF."
$pTestade.Range.Text = "Testade att köra med 100 lådor och 5 stycken skadade personer, det gick inte heller att köra med IPP så testade att köra med 100 lådor och 4 stycken personer."

# --- The bookmark paragraph (reached reliably via chained .Next(), since indices /
#     cached refs captured before these inserts would now point at the wrong spot)
#     gets a new sentence inserted right before the bookmark marks themselves. ---
$pBookmark = $pTestade.Next()
$bmRange = $pBookmark.Range
$bmRange.Collapse(1)
$bmRange.InsertBefore("Det gick att köra och tog ungefär lite mer än 2 minuter med IPP.")

# --- New paragraphs AFTER the bookmark paragraph ---
$pVar = New-ParaAfter $pBookmark
$pSlutsats = New-ParaAfter $pVar
$pDetTar = New-ParaAfter $pSlutsats
$pEmptyAfter = New-ParaAfter $pDetTar
$pLabb13 = New-ParaAfter $pEmptyAfter
$pEmptyA = New-ParaAfter $pLabb13
$pEmptyB = New-ParaAfter $pEmptyA
$pSpaces = New-ParaAfter $pEmptyB
$pEmptyC = New-ParaAfter $pSpaces

$pVar.Range.Text = "Var lite nyfiken på hur stor påverkan lådorna har så testade då att köra med 50 lådor och 5 personer vilket inte heller var körbart med IPP."
$pSlutsats.Range.Text = "Slutsatsen är att sökrymden ökar drastiskt och problemet får en hög fögreningsfaktor då vi ökar antalet sjuka personer i problemet. Det märktes avsevärt då vi ökade ett problem med 100 lådor och 3 sjuka person till 100 lådor och 4 sjuka personer, exekveringstiden gick ifrån 3 sekunder till mer än 2 minuter."
$pDetTar.Range.Text = "Det tar då mer tid att konstruera grafen och sökningen i A* som IPP använder sig av."
$pLabb13.Range.Text = "Labb1.3"
$pSpaces.Range.Text = "  "

# --- Apply bold to the two section headings last, so the formatting never
#     leaks forward into paragraphs created afterwards. ---
$pLabb12.Range.Bold = 1
$pLabb13.Range.Bold = 1
